$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "AVL.py",
    "AVLTree.java",
    "BestFirstSearch.java",
    "bin_search.py",
    "Binary Search.java",
    "Binary_Search (2).java",
    "binary_search (2).py",
    "binary_search (3).py",
    "binary_search (4).py",
    "binary_search (5).py",
    "binary_search (6).py",
    "binary_search (7).py",
    "binary_search (8).py",
    "binary_search (9).py",
    "binary_search_tree.py",
    "Binary_search.java",
    "binary_search.py",
    "binary-search-tree.js",
    "binary.py",
    "BinarySearch (2).java",
    "BinarySearch (2).js",
    "BinarySearch (2).py",
    "BinarySearch (3).java",
    "BinarySearch (3).js",
    "binarySearch (3).py",
    "BinarySearch (4).java",
    "binarySearch (4).js",
    "BinarySearch (5).java",
    "BinarySearch (5).js",
    "BinarySearch (6).java",
    "binarySearch (6).js",
    "BinarySearch (7).java",
    "binarySearch (7).js",
    "BinarySearch (8).java",
    "BinarySearch (9).java",
    "binarySearch (10).java",
    "BinarySearch (11).java"
)

$data = @(
    @(0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,1,0,0),
    @(0,0,0,0,0,0,1,0,0),
    @(0,0,0,0,0,0,1,0,0),
    @(0,0,0,0,0,0,1,0,0),
    @(0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,1,0,0),
    @(0,0,0,0,0,0,1,0,0),
    @(0,0,0,0,0,0,1,0,0),
    @(0,0,0,0,0,0,1,0,0),
    @(0,0,0,0,0,0,1,0,0),
    @(0,0,0,0,0,0,1,0,0),
    @(0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,1,0,0),
    @(0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,1,0,0),
    @(0,0,0,0,0,0,1,0,0),
    @(0,0,0,0,0,0,1,0,0),
    @(0,0,0,0,0,0,1,0,0),
    @(0,0,0,0,0,0,1,0,0),
    @(0,0,0,0,0,0,1,0,0),
    @(0,0,0,0,0,0,1,0,0),
    @(0,0,0,0,0,0,1,0,0),
    @(0,0,0,0,0,0,1,0,0),
    @(0,0,0,0,0,0,1,0,0),
    @(0,0,0,0,0,0,1,0,0),
    @(0,0,0,0,0,0,1,0,0),
    @(0,0,0,0,0,0,1,0,0),
    @(0,0,0,0,0,0,1,0,0),
    @(0,0,0,0,0,0,1,0,0),
    @(0,0,0,0,0,0,1,0,0),
    @(0,0,0,0,1,0,1,0,0),
    @(0,0,0,0,0,0,1,0,0)
)

$startRow = 365
for ($i = 0; $i -lt $names.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $names[$i]
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 2).Value = $row[$c]
    }
}

$ws.Activate()
$excel.ActiveWindow.ScrollRow = 385
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A402").Select()
